# Insert a new price-report row for "Feria Lagunitas de Puerto Montt" - Apio
# right before the current row 218 (A1:R342 -> A1:R343), pushing the old
# rows 218..342 down to 219..343.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 218..342 down by one, creating a blank row 218.
$ws.Rows("218:218").Insert()

# Populate the newly inserted row 218 with the new record.
$row = 218
$ws.Cells.Item($row, 1).Value2 = 4
$ws.Cells.Item($row, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value2 = "Los Lagos"
$ws.Cells.Item($row, 4).Value2 = 44873
$ws.Cells.Item($row, 5).Value2 = 10
$ws.Cells.Item($row, 6).Value2 = 100112017
$ws.Cells.Item($row, 7).Value2 = "Apio"
$ws.Cells.Item($row, 8).Value2 = "Americana (o)"
$ws.Cells.Item($row, 9).Value2 = "Primera"
$ws.Cells.Item($row, 10).Value2 = 40
$ws.Cells.Item($row, 11).Value2 = 13000
$ws.Cells.Item($row, 12).Value2 = 14000
$ws.Cells.Item($row, 13).Value2 = 13500
$ws.Cells.Item($row, 14).Value2 = "$/docena de matas"
$ws.Cells.Item($row, 15).Value2 = "Región de Coquimbo"
$ws.Cells.Item($row, 16).Value2 = 2250
$ws.Cells.Item($row, 17).Value2 = 6
$ws.Cells.Item($row, 18).Value2 = "Hortaliza"
